# Apply a row-shuffle to the "variable" columns (D, I, J, K, L, M, P) of the
# data rows (2..123) in the active worksheet. The fixed columns
# (A,B,C,E,F,G,H,N,O,Q,R) and the header row stay untouched. For each
# destination row, the values come from the source row indicated in
# $srcRows (1-based position corresponds to destination row - 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 123

# Mapping: destination row (index 0 -> row 2) gets its D/I/J/K/L/M/P values
# from this source row number.
$srcRows = @(11,59,20,60,85,117,19,80,74,62,64,121,53,76,102,112,78,108,44,89,13,99,61,100,81,3,58,25,120,98,17,109,42,68,29,30,47,23,65,9,105,118,115,66,40,107,71,73,70,46,97,116,28,96,90,111,69,16,101,4,67,94,122,86,24,15,37,38,10,48,34,104,32,5,8,57,43,84,22,36,87,91,123,6,103,113,26,56,18,114,52,63,31,110,54,72,27,119,79,92,41,35,45,93,21,51,95,82,49,14,75,55,106,50,83,39,88,77,33,7,12,2)

# Columns (1-based) that move together with each logical row.
$cols = @(4, 9, 10, 11, 12, 13, 16)

# Snapshot original values for the columns we are about to shuffle, keyed by
# row number, before we start overwriting anything.
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $original[$r] = $rowVals
}

# Now write the shuffled values back.
for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $srcRows[$i]
    $rowVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
}
